$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$blankStyle = $ws.Range("A1").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "313.48"
$ws.Range("D2").Style = $blankStyle
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.46%"
$ws.Range("E2").Style = $blankStyle
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2"
$ws.Range("G2").Style = $blankStyle

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.43"
$ws.Range("D3").Style = $blankStyle
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.70%"
$ws.Range("E3").Style = $blankStyle
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2"
$ws.Range("G3").Style = $blankStyle

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.075"
$ws.Range("D4").Style = $blankStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.79%"
$ws.Range("E4").Style = $blankStyle
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2"
$ws.Range("G4").Style = $blankStyle

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08048"
$ws.Range("D5").Style = $blankStyle
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.83%"
$ws.Range("E5").Style = $blankStyle
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "2"
$ws.Range("G5").Style = $blankStyle

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.180"
$ws.Range("D6").Style = $blankStyle
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.85%"
$ws.Range("E6").Style = $blankStyle
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2"
$ws.Range("G6").Style = $blankStyle

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.029"
$ws.Range("D7").Style = $blankStyle
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.25%"
$ws.Range("E7").Style = $blankStyle
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2"
$ws.Range("G7").Style = $blankStyle

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("B8").Style = $blankStyle
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C8").Style = $blankStyle
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9287"
$ws.Range("D8").Style = $blankStyle
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.56%"
$ws.Range("E8").Style = $blankStyle
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "2"
$ws.Range("G8").Style = $blankStyle

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("B9").Style = $blankStyle
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C9").Style = $blankStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1003"
$ws.Range("D9").Style = $blankStyle
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.78%"
$ws.Range("E9").Style = $blankStyle
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "2"
$ws.Range("G9").Style = $blankStyle

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("B10").Style = $blankStyle
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C10").Style = $blankStyle
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1877"
$ws.Range("D10").Style = $blankStyle
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.23%"
$ws.Range("E10").Style = $blankStyle
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "2"
$ws.Range("G10").Style = $blankStyle

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("B11").Style = $blankStyle
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C11").Style = $blankStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09127"
$ws.Range("D11").Style = $blankStyle
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.47%"
$ws.Range("E11").Style = $blankStyle
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "2"
$ws.Range("G11").Style = $blankStyle

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("B12").Style = $blankStyle
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C12").Style = $blankStyle
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03590"
$ws.Range("D12").Style = $blankStyle
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.75%"
$ws.Range("E12").Style = $blankStyle
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "2"
$ws.Range("G12").Style = $blankStyle

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("B13").Style = $blankStyle
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C13").Style = $blankStyle
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09936"
$ws.Range("D13").Style = $blankStyle
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.17%"
$ws.Range("E13").Style = $blankStyle
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2"
$ws.Range("G13").Style = $blankStyle

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("B14").Style = $blankStyle
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C14").Style = $blankStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001433"
$ws.Range("D14").Style = $blankStyle
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.20%"
$ws.Range("E14").Style = $blankStyle
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2"
$ws.Range("G14").Style = $blankStyle

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("B15").Style = $blankStyle
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").Style = $blankStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005656"
$ws.Range("D15").Style = $blankStyle
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.27%"
$ws.Range("E15").Style = $blankStyle
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "2"
$ws.Range("G15").Style = $blankStyle

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("B16").Style = $blankStyle
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C16").Style = $blankStyle
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.465"
$ws.Range("D16").Style = $blankStyle
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.18%"
$ws.Range("E16").Style = $blankStyle
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "2"
$ws.Range("G16").Style = $blankStyle

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("B17").Style = $blankStyle
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C17").Style = $blankStyle
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.137"
$ws.Range("D17").Style = $blankStyle
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.06%"
$ws.Range("E17").Style = $blankStyle
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "2"
$ws.Range("G17").Style = $blankStyle

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.798"
$ws.Range("D18").Style = $blankStyle
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.69%"
$ws.Range("E18").Style = $blankStyle
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "2"
$ws.Range("G18").Style = $blankStyle

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3375"
$ws.Range("D19").Style = $blankStyle
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.24%"
$ws.Range("E19").Style = $blankStyle
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "2"
$ws.Range("G19").Style = $blankStyle

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1377"
$ws.Range("D20").Style = $blankStyle
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.35%"
$ws.Range("E20").Style = $blankStyle
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "2"
$ws.Range("G20").Style = $blankStyle

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.059"
$ws.Range("D21").Style = $blankStyle
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.26%"
$ws.Range("E21").Style = $blankStyle
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "2"
$ws.Range("G21").Style = $blankStyle

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2332"
$ws.Range("D22").Style = $blankStyle
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.85%"
$ws.Range("E22").Style = $blankStyle
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "2"
$ws.Range("G22").Style = $blankStyle

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04592"
$ws.Range("D23").Style = $blankStyle
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.20%"
$ws.Range("E23").Style = $blankStyle
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "2"
$ws.Range("G23").Style = $blankStyle

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("D24").Style = $blankStyle
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.79%"
$ws.Range("E24").Style = $blankStyle
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "2"
$ws.Range("G24").Style = $blankStyle

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004758"
$ws.Range("D25").Style = $blankStyle
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.73%"
$ws.Range("E25").Style = $blankStyle
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "2"
$ws.Range("G25").Style = $blankStyle

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("D26").Style = $blankStyle
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.27%"
$ws.Range("E26").Style = $blankStyle
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "2"
$ws.Range("G26").Style = $blankStyle

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004494"
$ws.Range("D27").Style = $blankStyle
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "64.88%"
$ws.Range("E27").Style = $blankStyle
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "2"
$ws.Range("G27").Style = $blankStyle

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "2"
$ws.Range("G28").Style = $blankStyle

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "2"
$ws.Range("G29").Style = $blankStyle

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "2"
$ws.Range("G30").Style = $blankStyle

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "2"
$ws.Range("G31").Style = $blankStyle

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "2"
$ws.Range("G32").Style = $blankStyle

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "2"
$ws.Range("G33").Style = $blankStyle

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "2"
$ws.Range("G34").Style = $blankStyle

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "2"
$ws.Range("G35").Style = $blankStyle

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "2"
$ws.Range("G36").Style = $blankStyle

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "2"
$ws.Range("G37").Style = $blankStyle

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "2"
$ws.Range("G38").Style = $blankStyle

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01949"
$ws.Range("D39").Style = $blankStyle
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.53%"
$ws.Range("E39").Style = $blankStyle
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "2"
$ws.Range("G39").Style = $blankStyle

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04974"
$ws.Range("D40").Style = $blankStyle
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.59%"
$ws.Range("E40").Style = $blankStyle
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "2"
$ws.Range("G40").Style = $blankStyle

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007758"
$ws.Range("D41").Style = $blankStyle
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.39%"
$ws.Range("E41").Style = $blankStyle
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "2"
$ws.Range("G41").Style = $blankStyle

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1398"
$ws.Range("D42").Style = $blankStyle
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.21%"
$ws.Range("E42").Style = $blankStyle
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "2"
$ws.Range("G42").Style = $blankStyle

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007792"
$ws.Range("D43").Style = $blankStyle
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.98%"
$ws.Range("E43").Style = $blankStyle
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "2"
$ws.Range("G43").Style = $blankStyle

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("D44").Style = $blankStyle
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.62%"
$ws.Range("E44").Style = $blankStyle
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "2"
$ws.Range("G44").Style = $blankStyle

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01149"
$ws.Range("D45").Style = $blankStyle
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.41%"
$ws.Range("E45").Style = $blankStyle
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "2"
$ws.Range("G45").Style = $blankStyle

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006231"
$ws.Range("D46").Style = $blankStyle
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.74%"
$ws.Range("E46").Style = $blankStyle
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "2"
$ws.Range("G46").Style = $blankStyle

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("D47").Style = $blankStyle
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.16%"
$ws.Range("E47").Style = $blankStyle
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "2"
$ws.Range("G47").Style = $blankStyle

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.94"
$ws.Range("D48").Style = $blankStyle
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-24.96%"
$ws.Range("E48").Style = $blankStyle
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "2"
$ws.Range("G48").Style = $blankStyle

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001897"
$ws.Range("D49").Style = $blankStyle
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.14%"
$ws.Range("E49").Style = $blankStyle
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "2"
$ws.Range("G49").Style = $blankStyle

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("D50").Style = $blankStyle
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"
$ws.Range("E50").Style = $blankStyle
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "2"
$ws.Range("G50").Style = $blankStyle

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("D51").Style = $blankStyle
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.16%"
$ws.Range("E51").Style = $blankStyle
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "2"
$ws.Range("G51").Style = $blankStyle
